# Update JacobiDataSource test workbook to support formula as matrix dimension.
# - Rename the third sheet ("Sheet3") to "test Dynamic Dimension 1x7"
# - Populate it with an id cell, two COUNTA formulas, and a 1x7 data row
# - Make it the active sheet (moves tabSelected / activeTab from sheet 2 to sheet 3)

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "test Dynamic Dimension 1x7"

# Id cell
$ws3.Range("A1").Value = "#0"

# Formulas that compute the matrix dimension dynamically
$ws3.Range("A2").Formula = "=COUNTA(A3:A8)"
$ws3.Range("B2").Formula = "=COUNTA(3:3)"

# The 1x7 matrix row
$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = 2
$ws3.Range("C3").Value = 3
$ws3.Range("D3").Value = 4
$ws3.Range("E3").Value = 5
$ws3.Range("F3").Value = 6
$ws3.Range("G3").Value = 7

# Select B3 and make this the active sheet/tab
$ws3.Range("B3").Select()
$ws3.Activate()
